$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Change energy unit labels from kcal/mol to kJ/mol ---
$ws.Range("G2").Value = "Relative Energy / kJ/mol (Gas Phase)"
$ws.Range("H2").Value = "Relative Energy / kJ/mol (SMD)"
$ws.Range("O2").Value = "Relative Energy / kJ/mol (Gas Phase)"
$ws.Range("P2").Value = "Relative Energy / kJ/mol (PCM)"

# --- Update conversion-factor formulas: Hartree -> kcal/mol (627.5095) becomes Hartree -> kJ/mol (2625.5) ---
$ws.Range("G6").Formula = "=(E6-E3)*2625.5"
$ws.Range("H6").Formula = "=(D6-D3)*2625.5"
$ws.Range("G7").Formula = "=(E7-E4)*2625.5"
$ws.Range("G8").Formula = "=(E8-E5)*2625.5"
$ws.Range("O8").Formula = "=(M8-M5)*2625.5"
$ws.Range("P8").Formula = "=(N8-N5)*2625.5"

# --- Update view state: scroll/selection moved from D1/I5 to C1/I3 ---
$win = $excel.Windows.Item(1)
$ws.Range("I3").Select()
$win.ScrollColumn = 3
$win.ScrollRow = 1
